# Add a new "buy_condition_type" column before the existing "timezone" column
# (i.e. insert a new column X, pushing the old X/"timezone" column to Y),
# populate it with "1h" for every data row, and renumber the "id" values
# in column A for the last rows of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at X; this shifts the existing column X ("timezone")
# to column Y, matching the diff (dimension grows from A1:X39 to A1:Y39).
$ws.Columns("X:X").Insert()

# Header for the newly inserted column.
$ws.Range("X1").Value = "buy_condition_type"

# Populate the new column with "1h" for every data row (2 through 39).
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 24).Value = "1h"
}

# Renumber the trailing "id" values in column A (rows 32-39).
$ws.Range("A32").Value = 30
$ws.Range("A33").Value = 31
$ws.Range("A34").Value = 32
$ws.Range("A35").Value = 33
$ws.Range("A36").Value = 34
$ws.Range("A37").Value = 35
$ws.Range("A38").Value = 36
$ws.Range("A39").Value = 37
